# Update "Forecast Comparison" sheet with corrected forecast output:
#  - insert a new "Week_Start_Date" column (B), shifting ASIN..is_holiday_week right by one
#  - normalize the Week labels (drop the leading zero: W01 -> W1, ... W09 -> W9)
#  - populate the new Week_Start_Date column with the week's start date (kept as text)
#  - store is_holiday_week as a boolean value

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new column before column B (ASIN). This shifts B..I to C..J and
# keeps all existing values/formatting of those columns intact.
$ws.Columns.Item(2).Insert()

# New header for the inserted column.
$ws.Cells.Item(1, 2).Value = "Week_Start_Date"

$weeks = @("W1","W2","W3","W4","W5","W6","W7","W8","W9","W10","W11","W12","W13","W14","W15","W16")
$weekStartDates = @("2025-01-05","2025-01-12","2025-01-19","2025-01-26","2025-02-02","2025-02-09","2025-02-16","2025-02-23","2025-03-02","2025-03-09","2025-03-16","2025-03-23","2025-03-30","2025-04-06","2025-04-13","2025-04-20")

for ($i = 0; $i -lt $weeks.Count; $i++) {
    $row = $i + 2

    # Column A: week label without the leading zero.
    $ws.Cells.Item($row, 1).Value = $weeks[$i]

    # Column B: the new Week_Start_Date column. Force text formatting first so
    # the ISO-looking date string isn't auto-converted into a date serial.
    $cell = $ws.Cells.Item($row, 2)
    $cell.NumberFormat = "@"
    $cell.Value = $weekStartDates[$i]

    # Column J (was I): is_holiday_week, now stored as a boolean.
    $ws.Cells.Item($row, 10).Value = $false
}
